$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds an Excel date serial number for every data
# row (rows 2-265). The recorded value 45189 (2023-09-20) needs to be
# bumped to 45190 (2023-09-21) for all of them.
$lastRow = 265
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
